$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that flip from 0 to 1 per the diff
$targets = @("G3","H3","D4","E4","D5","E5","H6","H7","H8","D9","E9","H10","H11","H12","H13","H14","D15","E15","H16","D17","E17","H18")

foreach ($addr in $targets) {
    $ws.Range($addr).Value = 1
}
